# Co-optimization of H2 and Power system: add new Hydrogen-related plant
# technology rows (SR/SMR/SMR CCS/Electrolyzer/Fuel Cell/H2 Turbine) to the
# NewTechFramework sheet, rename the existing Hydrogen storage row's
# DataSource to "h2_storage", and apply a numeric format to H12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (plant type names) for the new rows 12-17 -------------------
$ws.Range("A12").Value = "SR"
$ws.Range("A13").Value = "SMR"
$ws.Range("A14").Value = "SMR CCS"
$ws.Range("A15").Value = "Electrolyzer"
$ws.Range("A16").Value = "Fuel Cell"
$ws.Range("A17").Value = "H2 Turbine"

# --- Column B (DataSource) ---------------------------------------------
# Order matters for shared-string layout: rows 16, 17, 15 first, then the
# existing Hydrogen row (11) is renamed, then the simple duplicates.
$ws.Range("B16").Value = "fuel_cell"
$ws.Range("B17").Value = "H2_turbine"
$ws.Range("B15").Value = "electrolyzer"
$ws.Range("B11").Value = "h2_storage"
$ws.Range("B12").Value = "SR"
$ws.Range("B13").Value = "SMR"
$ws.Range("B14").Value = "SMR"

# --- Column C (ATBTechnologyType) - all "NA" -----------------------------
$ws.Range("C12").Value = "NA"
$ws.Range("C13").Value = "NA"
$ws.Range("C14").Value = "NA"
$ws.Range("C15").Value = "NA"
$ws.Range("C16").Value = "NA"
$ws.Range("C17").Value = "NA"

# --- Column D (FuelType) -------------------------------------------------
$ws.Range("D12").Value = "Nuclear Fuel"
$ws.Range("D13").Value = "Natural Gas"
$ws.Range("D14").Value = "Natural Gas"
$ws.Range("D15").Value = "Electricity"
$ws.Range("D16").Value = "Hydrogen"
$ws.Range("D17").Value = "Hydrogen"

# --- Column E (ThermalOrRenewableOrStorage) ------------------------------
$ws.Range("E12").Value = "thermal"
$ws.Range("E13").Value = "h2"
$ws.Range("E14").Value = "h2"
$ws.Range("E15").Value = "h2"
$ws.Range("E16").Value = "h2toelectricity"
$ws.Range("E17").Value = "h2toelectricity"

# --- Column F (Capacity (MW)) --------------------------------------------
$ws.Range("F12").Value = 80
$ws.Range("F13").Value = 9170
$ws.Range("F14").Value = 9170
$ws.Range("F15").Value = 2000
$ws.Range("F16").Value = 50
$ws.Range("F17").Value = 240

# --- Column G (Heat Rate (Btu/kWh)) - only rows 15-17 --------------------
$ws.Range("G15").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("G17").Value = 0

# --- Column H - row 12 only gets a number format, no value --------------
$ws.Range("H12").NumberFormat = "#,##0.00"

# --- Column K (NSPSCompliant) --------------------------------------------
$ws.Range("K12").Value = "Yes"
$ws.Range("K13").Value = "Yes"
$ws.Range("K14").Value = "Yes"
$ws.Range("K15").Value = "Yes"
$ws.Range("K16").Value = "Yes"
$ws.Range("K17").Value = "Yes"

# --- Column L (NOxEmRate) -------------------------------------------------
$ws.Range("L12").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("L17").Value = 0

# --- Column M (SO2EmRate) -------------------------------------------------
$ws.Range("M12").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("M17").Value = 0

# --- Column N (CO2EmRate) -------------------------------------------------
$ws.Range("N12").Value = 0
$ws.Range("N13").Value = 20000
$ws.Range("N14").Value = 2000
$ws.Range("N15").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("N17").Value = 0

# --- Column O (Lifetime(years)) ------------------------------------------
$ws.Range("O12").Value = 40
$ws.Range("O13").Value = 25
$ws.Range("O14").Value = 25
$ws.Range("O15").Value = 10
$ws.Range("O16").Value = 10
$ws.Range("O17").Value = 25

# --- Column P (FuelPrice($/MMBtu)) - only rows 15-17 ---------------------
$ws.Range("P15").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("P17").Value = 0

# --- Column Q (SO2 Scrubber) ----------------------------------------------
$ws.Range("Q12").Value = "NA"
$ws.Range("Q13").Value = "NA"
$ws.Range("Q14").Value = "NA"
$ws.Range("Q15").Value = "NA"
$ws.Range("Q16").Value = "NA"
$ws.Range("Q17").Value = "NA"

# --- Column R (CoalType) ---------------------------------------------------
$ws.Range("R12").Value = "NA"
$ws.Range("R13").Value = "NA"
$ws.Range("R14").Value = "NA"
$ws.Range("R15").Value = "NA"
$ws.Range("R16").Value = "NA"
$ws.Range("R17").Value = "NA"

# --- Column S (Efficiency) -------------------------------------------------
$ws.Range("S12").Value = "NA"
$ws.Range("S13").Value = "NA"
$ws.Range("S14").Value = "NA"
$ws.Range("S15").Value = "NA"
$ws.Range("S16").Value = "NA"
$ws.Range("S17").Value = "NA"

# --- Column T (Nameplate Energy Capacity (MWh)) - only row 12 -------------
$ws.Range("T12").Value = "NA"

# --- Column U (Minimum Energy Capacity (MWh)) ------------------------------
$ws.Range("U12").Value = "NA"
$ws.Range("U13").Value = "NA"
$ws.Range("U14").Value = "NA"
$ws.Range("U15").Value = "NA"
$ws.Range("U16").Value = "NA"
$ws.Range("U17").Value = "NA"

# --- Column V (Maximum Charge Rate (MW)) -----------------------------------
$ws.Range("V12").Value = "NA"
$ws.Range("V13").Value = "NA"
$ws.Range("V14").Value = "NA"
$ws.Range("V15").Value = "NA"
$ws.Range("V16").Value = "NA"
$ws.Range("V17").Value = "NA"

# --- Column W (ECAPEX(2012$/MWH)) ------------------------------------------
$ws.Range("W12").Value = "NA"
$ws.Range("W13").Value = "NA"
$ws.Range("W14").Value = "NA"
$ws.Range("W15").Value = "NA"
$ws.Range("W16").Value = "NA"
$ws.Range("W17").Value = "NA"

# --- Column X (PlantCategory) ----------------------------------------------
$ws.Range("X12").Value = "Nuclear"
$ws.Range("X13").Value = "SMR"
$ws.Range("X14").Value = "SMR"
$ws.Range("X15").Value = "Electrolyzer"
$ws.Range("X16").Value = "FuelCell"
$ws.Range("X17").Value = "H2Turbine"

# --- Final selection / view state matching the saved workbook -------------
$ws.Range("I12").Select()
